$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string table entry that used to read "TekCollect" now reads
# "SmartCollect", the old "SmartCollect" entry is dropped, and a new
# "TekCollect" entry is appended at the end (after "Abellaregistration").
# Net visible effect on the sheet's App-Name column (B2:B5) is a one-row
# rotation of the four values.
$ws.Range("B2").Value = "SmartCollect"
$ws.Range("B3").Value = "CapitalAccounts"
$ws.Range("B4").Value = "Abellaregistration"
$ws.Range("B5").Value = "TekCollect"

# The header row's fill pattern changes from "lightTrellis" to "solid"
# (same indexed-46 fill colour as before).
$header = $ws.Range("A1:D1").Interior
$header.ColorIndex = 46
$header.Pattern = [Microsoft.Office.Interop.Excel.XlPatternType]::xlPatternSolid
